# Re-randomize the "conditions" trial table on Sheet1: several rows' stimulus
# color / duration / target / sound values were reshuffled, and nine new trial
# rows (43-51) were appended, extending the data range from A1:I42 to A1:I51.
# Only the cells that actually change are touched; the SUM(C:C) formula in I6
# and the sheet dimension recalculate automatically from the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Cells.Item(5,3).Value = 0.05
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = "silent.wav"
# Row 6
$ws.Cells.Item(6,2).Value = "red.png"
$ws.Cells.Item(6,4).Value = 1
$ws.Cells.Item(6,5).Value = 1
# Row 7
$ws.Cells.Item(7,2).Value = "orange.png"
$ws.Cells.Item(7,3).Value = 0.05
# Row 9
$ws.Cells.Item(9,2).Value = "green.png"
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = ""
$ws.Cells.Item(9,6).Value = "silent.wav"
# Row 11
$ws.Cells.Item(11,2).Value = "red.png"
$ws.Cells.Item(11,3).Value = 0.05
$ws.Cells.Item(11,5).Value = 0
# Row 12
$ws.Cells.Item(12,2).Value = "red.png"
$ws.Cells.Item(12,4).Value = 1
$ws.Cells.Item(12,5).Value = 1
# Row 13
$ws.Cells.Item(13,2).Value = "orange.png"
$ws.Cells.Item(13,3).Value = 0.05
# Row 14
$ws.Cells.Item(14,2).Value = "orange.png"
$ws.Cells.Item(14,6).Value = "beep.wav"
# Row 15
$ws.Cells.Item(15,2).Value = "green.png"
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = ""
$ws.Cells.Item(15,6).Value = "silent.wav"
# Row 17
$ws.Cells.Item(17,2).Value = "green.png"
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = ""
$ws.Cells.Item(17,6).Value = "silent.wav"
# Row 20
$ws.Cells.Item(20,2).Value = "orange.png"
$ws.Cells.Item(20,6).Value = "beep.wav"
# Row 22
$ws.Cells.Item(22,2).Value = "red.png"
$ws.Cells.Item(22,3).Value = 0.05
$ws.Cells.Item(22,5).Value = 0
# Row 23
$ws.Cells.Item(23,2).Value = "red.png"
$ws.Cells.Item(23,4).Value = 1
$ws.Cells.Item(23,5).Value = 1
$ws.Cells.Item(23,6).Value = "beep.wav"
# Row 24
$ws.Cells.Item(24,2).Value = "green.png"
$ws.Cells.Item(24,6).Value = "silent.wav"
# Row 25
$ws.Cells.Item(25,2).Value = "red.png"
$ws.Cells.Item(25,4).Value = 1
$ws.Cells.Item(25,5).Value = 1
$ws.Cells.Item(25,6).Value = "beep.wav"
# Row 27
$ws.Cells.Item(27,2).Value = "orange.png"
$ws.Cells.Item(27,6).Value = "beep.wav"
# Row 28
$ws.Cells.Item(28,2).Value = "green.png"
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = ""
$ws.Cells.Item(28,6).Value = "silent.wav"
# Row 29
$ws.Cells.Item(29,2).Value = "green.png"
$ws.Cells.Item(29,6).Value = "silent.wav"
# Row 30
$ws.Cells.Item(30,2).Value = "orange.png"
$ws.Cells.Item(30,6).Value = "beep.wav"
# Row 32
$ws.Cells.Item(32,2).Value = "green.png"
$ws.Cells.Item(32,6).Value = "silent.wav"
# Row 33
$ws.Cells.Item(33,2).Value = "green.png"
$ws.Cells.Item(33,4).Value = 0
$ws.Cells.Item(33,5).Value = ""
$ws.Cells.Item(33,6).Value = "silent.wav"
# Row 34
$ws.Cells.Item(34,2).Value = "red.png"
$ws.Cells.Item(34,4).Value = 1
$ws.Cells.Item(34,5).Value = 1
$ws.Cells.Item(34,6).Value = "beep.wav"
# Row 35
$ws.Cells.Item(35,2).Value = "orange.png"
$ws.Cells.Item(35,3).Value = 0.05
# Row 37
$ws.Cells.Item(37,2).Value = "green.png"
$ws.Cells.Item(37,6).Value = "silent.wav"
# Row 40
$ws.Cells.Item(40,2).Value = "red.png"
$ws.Cells.Item(40,3).Value = 0.05
$ws.Cells.Item(40,5).Value = 0
# Row 42
$ws.Cells.Item(42,3).Value = 3
# Row 43 (new)
$ws.Cells.Item(43,1).Value = 1
$ws.Cells.Item(43,2).Value = "green.png"
$ws.Cells.Item(43,3).Value = 3
$ws.Cells.Item(43,4).Value = 0
$ws.Cells.Item(43,6).Value = "silent.wav"
# Row 44 (new)
$ws.Cells.Item(44,1).Value = 1
$ws.Cells.Item(44,2).Value = "orange.png"
$ws.Cells.Item(44,3).Value = 3
$ws.Cells.Item(44,4).Value = 0
$ws.Cells.Item(44,6).Value = "beep.wav"
# Row 45 (new)
$ws.Cells.Item(45,1).Value = 1
$ws.Cells.Item(45,2).Value = "red.png"
$ws.Cells.Item(45,3).Value = 0.05
$ws.Cells.Item(45,4).Value = 0
$ws.Cells.Item(45,6).Value = "silent.wav"
# Row 46 (new)
$ws.Cells.Item(46,1).Value = 1
$ws.Cells.Item(46,2).Value = "red.png"
$ws.Cells.Item(46,3).Value = 3
$ws.Cells.Item(46,4).Value = 0
$ws.Cells.Item(46,6).Value = "beep.wav"
# Row 47 (new)
$ws.Cells.Item(47,1).Value = 1
$ws.Cells.Item(47,2).Value = "green.png"
$ws.Cells.Item(47,3).Value = 3
$ws.Cells.Item(47,4).Value = 0
$ws.Cells.Item(47,6).Value = "silent.wav"
# Row 48 (new)
$ws.Cells.Item(48,1).Value = 1
$ws.Cells.Item(48,2).Value = "orange.png"
$ws.Cells.Item(48,3).Value = 3
$ws.Cells.Item(48,4).Value = 0
$ws.Cells.Item(48,6).Value = "beep.wav"
# Row 49 (new)
$ws.Cells.Item(49,1).Value = 1
$ws.Cells.Item(49,2).Value = "green.png"
$ws.Cells.Item(49,3).Value = 3
$ws.Cells.Item(49,4).Value = 0
$ws.Cells.Item(49,6).Value = "silent.wav"
# Row 50 (new)
$ws.Cells.Item(50,1).Value = 1
$ws.Cells.Item(50,2).Value = "red.png"
$ws.Cells.Item(50,3).Value = 3
$ws.Cells.Item(50,4).Value = 1
$ws.Cells.Item(50,5).Value = 1
$ws.Cells.Item(50,6).Value = "beep.wav"
# Row 51 (new)
$ws.Cells.Item(51,1).Value = 1
$ws.Cells.Item(51,2).Value = "green.png"
$ws.Cells.Item(51,3).Value = 2.9
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,6).Value = "silent.wav"

# Match the saved selection state (A2:F51) recorded in the target workbook.
$ws.Range("A2:F51").Select() | Out-Null
